$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new board/project name.
$ws.Name = "BOM_ESP32-Wiegand-Hat_ESP32-Wie"

# "Board" designator text is now uppercase.
$ws.Range("D3").Value = "BOARD"

# The DIP switch (item 4, DP-02RP) used to have a single BOM row labeled
# "DBG". It is now split into two rows describing the switch's two
# positions: "ON" (keeps the original row's footprint/value columns only)
# and a brand-new "OFF" row carrying the remaining BOM columns
# (footprint, manufacturer, supplier part, supplier).
$ws.Range("D5").Value = "ON"

# Insert a new row right after row 5 for the "OFF" entry; this shifts all
# subsequent rows (old rows 6-10) down by one.
$ws.Rows.Item(6).Insert()

# Clear any stale values that might have been duplicated into row 5 columns
# E:J by the insert/shift (row 5 should only have columns A:D populated).
$ws.Range("E5:J5").ClearContents()

# Populate the new "OFF" row (row 6) with columns A:G.
$ws.Range("A6").Value = "OFF"
$ws.Range("B6").Value = "SW-TH_DP-02XP"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "DP-02RP"
$ws.Range("E6").Value = "韩国韩荣"
$ws.Range("F6").Value = "C129041"
$ws.Range("G6").Value = "LCSC"

# "Wiegand" designator text is now uppercase (this row shifted from 9 to 10).
$ws.Range("D10").Value = "WIEGAND"
